# Auto-applies the diff for portugal_liga-3_2023-2024.xlsx
# - Reorders several match rows within existing matchday blocks (rows 27-29, 44-50, 94-98)
#   so that their F:V (match details) content matches the updated source order.
# - Appends a brand-new match row (row 100) for Trofense vs Felgueiras (08-10/11/2023).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27
$ws.Range("F27").Value = "Braga B"
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = "Trofense"
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 2.84
$ws.Range("K27").Value = "19/08/2023 03:43"
$ws.Range("L27").Value = 1.95
$ws.Range("M27").Value = "20/08/2023 18:51"
$ws.Range("N27").Value = 3.09
$ws.Range("O27").Value = "19/08/2023 03:43"
$ws.Range("P27").Value = 3.41
$ws.Range("Q27").Value = "20/08/2023 18:51"
$ws.Range("R27").Value = 2.64
$ws.Range("S27").Value = "19/08/2023 03:43"
$ws.Range("T27").Value = 4.17
$ws.Range("U27").Value = "20/08/2023 18:51"
$ws.Range("V27").Value = "https://www.betexplorer.com/football/portugal/liga-3/braga-trofense/0SxpaNx9/"

# Row 28
$ws.Range("F28").Value = "Canelas 2010"
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = "Felgueiras"
$ws.Range("I28").Value = 8
$ws.Range("J28").Value = 2.09
$ws.Range("K28").Value = "19/08/2023 03:43"
$ws.Range("L28").Value = 2.67
$ws.Range("M28").Value = "20/08/2023 18:51"
$ws.Range("N28").Value = 3.24
$ws.Range("O28").Value = "19/08/2023 03:43"
$ws.Range("P28").Value = 3.54
$ws.Range("Q28").Value = "20/08/2023 12:29"
$ws.Range("R28").Value = 3.52
$ws.Range("S28").Value = "19/08/2023 03:43"
$ws.Range("T28").Value = 2.59
$ws.Range("U28").Value = "20/08/2023 18:51"
$ws.Range("V28").Value = "https://www.betexplorer.com/football/portugal/liga-3/canelas-2010-fc-felgueiras/QgUi3ZYI/"

# Row 29
$ws.Range("F29").Value = "Anadia"
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = "Sanjoanense"
$ws.Range("I29").Value = 2
$ws.Range("J29").Value = 1.91
$ws.Range("K29").Value = "19/08/2023 03:43"
$ws.Range("L29").Value = 1.75
$ws.Range("M29").Value = "20/08/2023 18:53"
$ws.Range("N29").Value = 3.35
$ws.Range("O29").Value = "19/08/2023 03:43"
$ws.Range("P29").Value = 3.85
$ws.Range("Q29").Value = "20/08/2023 18:53"
$ws.Range("R29").Value = 3.98
$ws.Range("S29").Value = "19/08/2023 03:43"
$ws.Range("T29").Value = 4.63
$ws.Range("U29").Value = "20/08/2023 18:53"
$ws.Range("V29").Value = "https://www.betexplorer.com/football/portugal/liga-3/anadia-sanjoanense/4vtq5e46/"

# Row 44
$ws.Range("F44").Value = "Covilha"
$ws.Range("G44").Value = 1
$ws.Range("H44").Value = "Alverca"
$ws.Range("I44").Value = 2
$ws.Range("J44").Value = 1.79
$ws.Range("K44").Value = "30/08/2023 16:13"
$ws.Range("L44").Value = 1.93
$ws.Range("M44").Value = "03/09/2023 17:51"
$ws.Range("N44").Value = 3.44
$ws.Range("O44").Value = "30/08/2023 16:13"
$ws.Range("P44").Value = 3.47
$ws.Range("Q44").Value = "03/09/2023 17:51"
$ws.Range("R44").Value = 4.81
$ws.Range("S44").Value = "30/08/2023 16:13"
$ws.Range("T44").Value = 4.14
$ws.Range("U44").Value = "03/09/2023 17:51"
$ws.Range("V44").Value = "https://www.betexplorer.com/football/portugal/liga-3/covilha-alverca/jVMIITjd/"

# Row 45
$ws.Range("F45").Value = "Pero Pinheiro"
$ws.Range("G45").Value = 2
$ws.Range("H45").Value = "Amora"
$ws.Range("I45").Value = 3
$ws.Range("J45").Value = 2.86
$ws.Range("K45").Value = "30/08/2023 16:13"
$ws.Range("L45").Value = 3.75
$ws.Range("M45").Value = "03/09/2023 17:22"
$ws.Range("N45").Value = 3.08
$ws.Range("O45").Value = "30/08/2023 16:13"
$ws.Range("P45").Value = 3.24
$ws.Range("Q45").Value = "03/09/2023 17:22"
$ws.Range("R45").Value = 2.53
$ws.Range("S45").Value = "30/08/2023 16:13"
$ws.Range("T45").Value = 2.13
$ws.Range("U45").Value = "03/09/2023 17:22"
$ws.Range("V45").Value = "https://www.betexplorer.com/football/portugal/liga-3/pero-pinheiro-amora/xSIEJ9yj/"

# Row 46
$ws.Range("F46").Value = "Sporting CP B"
$ws.Range("G46").Value = 1
$ws.Range("H46").Value = "Academica"
$ws.Range("I46").Value = 2
$ws.Range("J46").Value = 2
$ws.Range("K46").Value = "30/08/2023 16:13"
$ws.Range("L46").Value = 1.8
$ws.Range("M46").Value = "03/09/2023 17:47"
$ws.Range("N46").Value = 3.33
$ws.Range("O46").Value = "30/08/2023 16:13"
$ws.Range("P46").Value = 3.74
$ws.Range("Q46").Value = "03/09/2023 17:42"
$ws.Range("R46").Value = 3.94
$ws.Range("S46").Value = "30/08/2023 16:13"
$ws.Range("T46").Value = 4.47
$ws.Range("U46").Value = "03/09/2023 17:47"
$ws.Range("V46").Value = "https://www.betexplorer.com/football/portugal/liga-3/sporting-lisbon-academica/bJSvQBbM/"

# Row 47
$ws.Range("F47").Value = "Caldas"
$ws.Range("G47").Value = 1
$ws.Range("H47").Value = "Oliveira Hospital"
$ws.Range("I47").Value = 1
$ws.Range("J47").Value = 1.95
$ws.Range("K47").Value = "30/08/2023 16:13"
$ws.Range("L47").Value = 1.78
$ws.Range("M47").Value = "03/09/2023 13:30"
$ws.Range("N47").Value = 3.3
$ws.Range("O47").Value = "30/08/2023 16:13"
$ws.Range("P47").Value = 3.6
$ws.Range("Q47").Value = "03/09/2023 13:30"
$ws.Range("R47").Value = 4.15
$ws.Range("S47").Value = "30/08/2023 16:13"
$ws.Range("T47").Value = 4.8
$ws.Range("U47").Value = "03/09/2023 13:30"
$ws.Range("V47").Value = "https://www.betexplorer.com/football/portugal/liga-3/caldas-sc-oliveira-hospital/hfsrPVDS/"

# Row 48
$ws.Range("F48").Value = "Canelas 2010"
$ws.Range("G48").Value = 2
$ws.Range("H48").Value = "SC Vianense"
$ws.Range("I48").Value = 1
$ws.Range("J48").Value = 1.9
$ws.Range("K48").Value = "30/08/2023 16:13"
$ws.Range("L48").Value = 2.14
$ws.Range("M48").Value = "03/09/2023 17:56"
$ws.Range("N48").Value = 3.37
$ws.Range("O48").Value = "30/08/2023 16:13"
$ws.Range("P48").Value = 3.41
$ws.Range("Q48").Value = "03/09/2023 17:56"
$ws.Range("R48").Value = 4.27
$ws.Range("S48").Value = "30/08/2023 16:13"
$ws.Range("T48").Value = 3.52
$ws.Range("U48").Value = "03/09/2023 17:56"
$ws.Range("V48").Value = "https://www.betexplorer.com/football/portugal/liga-3/canelas-2010-sc-vianense/bNGNGh3I/"

# Row 50
$ws.Range("F50").Value = "AD Fafe"
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = "Lusitania FC"
$ws.Range("I50").Value = 3
$ws.Range("J50").Value = 3.16
$ws.Range("K50").Value = "30/08/2023 16:13"
$ws.Range("L50").Value = 3.18
$ws.Range("M50").Value = "03/09/2023 17:34"
$ws.Range("N50").Value = 3.02
$ws.Range("O50").Value = "30/08/2023 16:13"
$ws.Range("P50").Value = 3.27
$ws.Range("Q50").Value = "03/09/2023 17:34"
$ws.Range("R50").Value = 2.46
$ws.Range("S50").Value = "30/08/2023 16:13"
$ws.Range("T50").Value = 2.36
$ws.Range("U50").Value = "03/09/2023 17:34"
$ws.Range("V50").Value = "https://www.betexplorer.com/football/portugal/liga-3/ad-fafe-lusitania-fc/pb1IJfJa/"

# Row 94
$ws.Range("F94").Value = "Sanjoanense"
$ws.Range("G94").Value = 3
$ws.Range("H94").Value = "Trofense"
$ws.Range("I94").Value = 2
$ws.Range("J94").Value = 2.67
$ws.Range("K94").Value = "04/11/2023 18:01"
$ws.Range("L94").Value = 3.01
$ws.Range("M94").Value = "05/11/2023 15:59"
$ws.Range("N94").Value = 3.07
$ws.Range("O94").Value = "04/11/2023 18:01"
$ws.Range("P94").Value = 3.18
$ws.Range("Q94").Value = "05/11/2023 15:59"
$ws.Range("R94").Value = 2.71
$ws.Range("S94").Value = "04/11/2023 18:01"
$ws.Range("T94").Value = 2.52
$ws.Range("U94").Value = "05/11/2023 15:59"
$ws.Range("V94").Value = "https://www.betexplorer.com/football/portugal/liga-3/sanjoanense-trofense/2X6FqrwM/"

# Row 95
$ws.Range("F95").Value = "1º Dezembro"
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = "Academica"
$ws.Range("I95").Value = 1
$ws.Range("J95").Value = 3.72
$ws.Range("K95").Value = "04/11/2023 18:02"
$ws.Range("L95").Value = 4.39
$ws.Range("M95").Value = "05/11/2023 15:53"
$ws.Range("N95").Value = 3.37
$ws.Range("O95").Value = "04/11/2023 18:02"
$ws.Range("P95").Value = 3.69
$ws.Range("Q95").Value = "05/11/2023 15:53"
$ws.Range("R95").Value = 1.97
$ws.Range("S95").Value = "04/11/2023 18:02"
$ws.Range("T95").Value = 1.83
$ws.Range("U95").Value = "05/11/2023 15:53"
$ws.Range("V95").Value = "https://www.betexplorer.com/football/portugal/liga-3/1-dezembro-academica/pCSjlLts/"

# Row 96
$ws.Range("F96").Value = "Atletico CP"
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = "Caldas"
$ws.Range("I96").Value = 2
$ws.Range("J96").Value = 1.95
$ws.Range("K96").Value = "04/11/2023 18:02"
$ws.Range("L96").Value = 2.15
$ws.Range("M96").Value = "05/11/2023 15:59"
$ws.Range("N96").Value = 3.38
$ws.Range("O96").Value = "04/11/2023 18:02"
$ws.Range("P96").Value = 3.34
$ws.Range("Q96").Value = "05/11/2023 15:59"
$ws.Range("R96").Value = 3.78
$ws.Range("S96").Value = "04/11/2023 18:02"
$ws.Range("T96").Value = 3.57
$ws.Range("U96").Value = "05/11/2023 15:59"
$ws.Range("V96").Value = "https://www.betexplorer.com/football/portugal/liga-3/atletico-cp-caldas-sc/xKpbnaBg/"

# Row 97
$ws.Range("F97").Value = "Covilha"
$ws.Range("G97").Value = 2
$ws.Range("H97").Value = "Pero Pinheiro"
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 1.54
$ws.Range("K97").Value = "04/11/2023 18:02"
$ws.Range("L97").Value = 1.4
$ws.Range("M97").Value = "05/11/2023 15:51"
$ws.Range("N97").Value = 4.12
$ws.Range("O97").Value = "04/11/2023 18:02"
$ws.Range("P97").Value = 4.81
$ws.Range("Q97").Value = "05/11/2023 15:51"
$ws.Range("R97").Value = 5.34
$ws.Range("S97").Value = "04/11/2023 18:02"
$ws.Range("T97").Value = 7.75
$ws.Range("U97").Value = "05/11/2023 15:51"
$ws.Range("V97").Value = "https://www.betexplorer.com/football/portugal/liga-3/covilha-pero-pinheiro/rcofmudm/"

# Row 98
$ws.Range("F98").Value = "Anadia"
$ws.Range("G98").Value = 1
$ws.Range("H98").Value = "Canelas 2010"
$ws.Range("I98").Value = 4
$ws.Range("J98").Value = 2.15
$ws.Range("K98").Value = "04/11/2023 18:02"
$ws.Range("L98").Value = 2.26
$ws.Range("M98").Value = "05/11/2023 15:59"
$ws.Range("N98").Value = 3.42
$ws.Range("O98").Value = "04/11/2023 18:02"
$ws.Range("P98").Value = 3.42
$ws.Range("Q98").Value = "05/11/2023 15:59"
$ws.Range("R98").Value = 3.17
$ws.Range("S98").Value = "04/11/2023 18:02"
$ws.Range("T98").Value = 3.23
$ws.Range("U98").Value = "05/11/2023 15:59"
$ws.Range("V98").Value = "https://www.betexplorer.com/football/portugal/liga-3/anadia-canelas-2010/K8HApONF/"

# Add new row 100 (new match result appended at the end of the sheet)
$ws.Range("A99:V99").Copy($ws.Range("A100:V100"))

$ws.Range("A100").Value = 99
$ws.Range("B100").Value = "portugal"
$ws.Range("C100").Value = "liga-3"
$ws.Range("D100").Value = "2023-2024"
$ws.Range("E100").Value = 45240.85416666666
$ws.Range("F100").Value = "Trofense"
$ws.Range("G100").Value = 2
$ws.Range("H100").Value = "Felgueiras"
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 3.86
$ws.Range("K100").Value = "08/11/2023 14:42"
$ws.Range("L100").Value = 4.44
$ws.Range("M100").Value = "10/11/2023 19:32"
$ws.Range("N100").Value = 3.34
$ws.Range("O100").Value = "08/11/2023 14:42"
$ws.Range("P100").Value = 3.65
$ws.Range("Q100").Value = "10/11/2023 19:32"
$ws.Range("R100").Value = 2.01
$ws.Range("S100").Value = "08/11/2023 14:42"
$ws.Range("T100").Value = 1.83
$ws.Range("U100").Value = "10/11/2023 15:03"
$ws.Range("V100").Value = "https://www.betexplorer.com/football/portugal/liga-3/trofense-fc-felgueiras/YFcovM0q/"
